# Reformatting of df_rejections: reorder/rename columns and insert a new
# "RESEAU" column, per the author's commit message
# ("reformating of df_rejections as asked").
#
# Target layout (A:H): FILIALE, RESEAU, ARN, Autorisation, Date Transaction,
# Montant, Devise, Motif
#
# Original layout (A:G): Montant, Devise, FILIALE, Transaction Date, ARN,
# Autorisation, Description

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Step 1: grow the table to the new A1:H4 footprint (must happen before
# the header text is written, so the new H column picks up a name from
# its header cell automatically).
# ---------------------------------------------------------------------
$lo.Resize($ws.Range("A1:H4"))

# ---------------------------------------------------------------------
# Step 2: write the header row (table column names sync automatically
# from the header cell text).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "FILIALE"
$ws.Range("B1").Value = "RESEAU"
$ws.Range("C1").Value = "ARN"
$ws.Range("D1").Value = "Autorisation"
$ws.Range("E1").Value = "Date Transaction"
$ws.Range("F1").Value = "Montant"
$ws.Range("G1").Value = "Devise"
$ws.Range("H1").Value = "Motif"

# ---------------------------------------------------------------------
# Step 3: write the body rows. Values that look numeric/date-like get a
# leading apostrophe so Excel keeps them as literal text (preserves
# leading zeros / long digit strings / the "yyyy-mm-dd" text). The
# formatting pass below (step 4) runs last and normalises the cell
# style, clearing any transient quote-prefix flag this introduces.
# ---------------------------------------------------------------------
$rows = @(
    @{ Filiale = "SG - COTE D IVOIRE"; Reseau = "MASTERCARD INTERNATIONAL"; Arn = "72681594150101332418418"; Autorisation = "059369"; Date = "2024-05-28"; Montant = "84,000"; Devise = "XOF"; Motif = "COUNTRY CODE INVALID FOR BUSINESS SERVICE ARRANGEMENT 4384001 AND ACCEPTANCE BRAND ID CODE 00000118 D0043 S06 DMC. INTERCHANGE RATE DESIGNATOR AND PROCESSING CODE/REVERSAL INDICATOR COMBINATION INVALID FOR 00000118 P0158 S04 BUSINESS SERVICE ARRANGEMENT 2060001 AND ACCEPTANCE BRAND ID CODE DMC." },
    @{ Filiale = "SG - COTE D IVOIRE"; Reseau = "MASTERCARD INTERNATIONAL"; Arn = "72681594150101332421271"; Autorisation = "059347"; Date = "2024-05-28"; Montant = "25,000"; Devise = "XOF"; Motif = "COUNTRY CODE INVALID FOR BUSINESS SERVICE ARRANGEMENT 4384001 AND ACCEPTANCE BRAND ID CODE 00000124 D0043 S06 DMC. INTERCHANGE RATE DESIGNATOR AND PROCESSING CODE/REVERSAL INDICATOR COMBINATION INVALID FOR 00000124 P0158 S04 BUSINESS SERVICE ARRANGEMENT 2060001 AND ACCEPTANCE BRAND ID CODE DMC." },
    @{ Filiale = "SG - COTE D IVOIRE"; Reseau = "MASTERCARD INTERNATIONAL"; Arn = "72681594150101332383190"; Autorisation = "059403"; Date = "2024-05-28"; Montant = "435,000"; Devise = "XOF"; Motif = "COUNTRY CODE INVALID FOR BUSINESS SERVICE ARRANGEMENT 4384001 AND ACCEPTANCE BRAND ID CODE 00000258 D0043 S06 DMC. INTERCHANGE RATE DESIGNATOR AND PROCESSING CODE/REVERSAL INDICATOR COMBINATION INVALID FOR 00000258 P0158 S04 BUSINESS SERVICE ARRANGEMENT 2060001 AND ACCEPTANCE BRAND ID CODE DMC." }
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.Filiale
    $ws.Range("B$r").Value = $row.Reseau
    $ws.Range("C$r").Value = "'" + $row.Arn
    $ws.Range("D$r").Value = "'" + $row.Autorisation
    $ws.Range("E$r").Value = "'" + $row.Date
    $ws.Range("F$r").Value = "'" + $row.Montant
    $ws.Range("G$r").Value = $row.Devise
    $ws.Range("H$r").Value = $row.Motif
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Step 4: re-apply formatting now that every cell holds its final value.
# Montant's header/value look moves from column A to column F; every
# other header/body cell gets the "generic" look that column B already
# carries. Copying A1/A2:A4 to F first (while A still looks right)
# avoids clobbering the source before it has been read.
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial($xlPasteFormats)

foreach ($r in 2..4) {
    $ws.Range("A$r").Copy()
    $ws.Range("F$r").PasteSpecial($xlPasteFormats)
}

foreach ($col in @("A", "C", "D", "E", "G", "H")) {
    $ws.Range("B1").Copy()
    $ws.Range("$col`1").PasteSpecial($xlPasteFormats)
}

foreach ($r in 2..4) {
    foreach ($col in @("A", "C", "D", "E", "G", "H")) {
        $ws.Range("B$r").Copy()
        $ws.Range("$col$r").PasteSpecial($xlPasteFormats)
    }
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 5: column widths for the new A:H layout.
# ---------------------------------------------------------------------
$widths = @(20, 26, 25, 14, 18, 9, 8, 297)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - (5 / 6)
}
